# Add a "cfop" column to the "PI hours" sheet and a new "cfop hours" sheet,
# matching the source-repo commit that introduced CFOP tracking.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$piSheet = $wb.Worksheets.Item("PI hours")

# A cell that already carries the bold/centered/bordered "header" style.
$headerStyleSrc = $piSheet.Cells.Item(1, 6)   # F1 ("app" header)
# A cell that already carries the bold/centered/bordered "index" style
# used on the leftmost numeric index column (A2:A7 etc.).
$indexStyleSrc = $piSheet.Cells.Item(2, 1)    # A2

# ---------------------------------------------------------------------------
# 1. "PI hours" sheet: add a new "cfop" column (G) with a per-PI cfop tag.
# ---------------------------------------------------------------------------
$piHeaderCell = $piSheet.Cells.Item(1, 7)     # G1
$headerStyleSrc.Copy()
$piHeaderCell.PasteSpecial($xlPasteFormats)
$piHeaderCell.Value = "cfop"

# Per-row cfop values (rows 2..7 correspond to the six PIs already listed)
$cfopValues = @(
    "['cfop_NH']",
    "['cfop_PARK']",
    "['cfop_WISSA']",
    "['cfop_SELIG']",
    "['cfop_MITRA']",
    "['cfop_GC']"
)

for ($i = 0; $i -lt $cfopValues.Length; $i++) {
    $row = $i + 2
    $piSheet.Cells.Item($row, 7).Value = $cfopValues[$i]
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. New "cfop hours" sheet, appended after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cfopSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$cfopSheet.Name = "cfop hours"

# Header row (B1:D1 -> cfop / hours / percentage)
$headers = @("cfop", "hours", "percentage")
for ($col = 2; $col -le 4; $col++) {
    $cell = $cfopSheet.Cells.Item(1, $col)
    $headerStyleSrc.Copy()
    $cell.PasteSpecial($xlPasteFormats)
    $cell.Value = $headers[$col - 2]
}

$excel.CutCopyMode = $false

# Data rows: index, cfop name, hours, percentage
$rows = @(
    @(0, "cfop_NH",    45,   57.69230769230769),
    @(1, "cfop_PARK",  9.5,  12.17948717948718),
    @(2, "cfop_WISSA", 8,    10.25641025641026),
    @(3, "cfop_MITRA", 6,    7.692307692307693),
    @(4, "cfop_GC",    5,    6.410256410256411),
    @(5, "cfop_SELIG", 4.5,  5.769230769230769)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $i + 2
    $data = $rows[$i]

    $indexCell = $cfopSheet.Cells.Item($row, 1)
    $indexStyleSrc.Copy()
    $indexCell.PasteSpecial($xlPasteFormats)
    $indexCell.Value = $data[0]

    $cfopSheet.Cells.Item($row, 2).Value = $data[1]
    $cfopSheet.Cells.Item($row, 3).Value = $data[2]
    $cfopSheet.Cells.Item($row, 4).Value = $data[3]
}

$excel.CutCopyMode = $false

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("PI hours").Activate()
